$d = $word.ActiveDocument

# --- Change 1 & 2: insert a centered bold heading paragraph before the
# existing first paragraph, and split the run around "methods" with
# proofErr gramStart/gramEnd markers (both land inside the first
# paragraph's Range, so they are done in a single InsertXML call that
# supplies two replacement paragraphs).
$p1 = $d.Paragraphs(1).Range

$frag1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>params ref and out modifiers</w:t></w:r></w:p>' + `
  '<w:p>' + `
    '<w:r><w:t xml:space="preserve"> params ref and out modifiers</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">, </w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">params </w:t></w:r>' + `
    '<w:r><w:t>fall</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:t>into</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> modifiers fall into modifier category when we talk about the </w:t></w:r>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:t>methods</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> and we must already know what' + [char]0x2019 + 's what are the methods I' + [char]0x2019 + 'd like functions in structure language there are methods in C sharp</w:t></w:r>' + `
    '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$p1.InsertXML($frag1)

# --- Change 3: append a new trailing run to the "Ref is used ..."
# paragraph (dropping the old trailing period and adding new text).
$rng = $d.Content
$rng.Find.Execute("Ref is used", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p2 = $rng.Paragraphs(1).Range

$apos = [char]0x2019
$part1 = "Ref is used  when we want to pass a variable with its reference, I find it weird tbh. Since I guess it" + $apos + "s okay in the structured language but in OOp we have many other better way to do that. I am kinda sure we don" + $apos + "t need in our C# practice but good to know"
$part2 = ", since if someone has used it in our project already."

$frag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
  '<w:p><w:r><w:t>' + $part1 + '</w:t></w:r><w:r><w:t>' + $part2 + '</w:t></w:r></w:p>' + `
  '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$p2.InsertXML($frag2)

Write-Host "Edit complete"
